$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.175.47"
$ws.Range("E2").Value = "  -2.36%  "

$ws.Range("D3").Value = "3.478.89"
$ws.Range("E3").Value = "  -4.15%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.82%  "

$ws.Range("D7").Value = "3.475.83"
$ws.Range("E7").Value = "  -4.30%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.40%  "

$ws.Range("D14").Value = "4.065.79"
$ws.Range("E14").Value = "  -3.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.86%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.187.47"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.465.57"
$ws.Range("E17").Value = "  -4.66%  "

$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.613"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "3.616.05"
$ws.Range("E26").Value = "  -4.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.14%  "

$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.161"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "

$ws.Range("D35").Value = "3.463.20"
$ws.Range("E35").Value = "  -4.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "177.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.49%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0868"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.869"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.974"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.38%  "
